$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
$wsALC.Range("H41").Value = 880.9231
$wsALC.Range("I41").Value = 820.3333
$wsALC.Range("K41").Value = 820.3333
$wsALC.Range("M41").Value = -380.3333
$wsALC.Range("H69").Value = 90916270
$wsALC.Range("J69").Value = 90916270
$wsALC.Range("L69").Value = 272748810
$wsALC.Range("N69").Value = -272750558
$wsALC.Range("H72").Value = 90916270
$wsALC.Range("J72").Value = 90916270
$wsALC.Range("L72").Value = 818246430
$wsALC.Range("N72").Value = -818255166
$wsALC.Range("H107").Value = 16862
$wsALC.Range("I107").Value = 17748.2
$wsALC.Range("K107").Value = 17748.2
$wsALC.Range("M107").Value = -15828.2
$wsALC.Range("H113").Value = 19875
$wsALC.Range("I113").Value = 60000
$wsALC.Range("K113").Value = 60000
$wsALC.Range("M113").Value = -56746
$wsALC.Range("H132").Value = 1993.2106
$wsALC.Range("I132").Value = 1591.5625
$wsALC.Range("K132").Value = 4774.6875
$wsALC.Range("M132").Value = -2244.6875
$wsALC.Range("H137").Value = 5074.4243
$wsALC.Range("I137").Value = 6206.56
$wsALC.Range("K137").Value = 18619.68
$wsALC.Range("M137").Value = -16069.68

# --- ARM ---
$wsARM.Range("H2").Value = 1544.2693
$wsARM.Range("I2").Value = 571.35297
$wsARM.Range("K2").Value = 571.35297
$wsARM.Range("M2").Value = -458.35297
$wsARM.Range("H32").Value = 2324.7407
$wsARM.Range("I32").Value = 2274.2642
$wsARM.Range("K32").Value = 2274.2642
$wsARM.Range("M32").Value = -1987.2642
$wsARM.Range("H61").Value = 3269
$wsARM.Range("I61").Value = 3269
$wsARM.Range("J61").Value = 0
$wsARM.Range("K61").Value = 3269
$wsARM.Range("L61").Value = 0
$wsARM.Range("M61").Value = -3057
$wsARM.Range("N61").ClearContents()
$wsARM.Range("H63").Value = 4749.8
$wsARM.Range("I63").Value = 775
$wsARM.Range("K63").Value = 775
$wsARM.Range("M63").Value = -89
$wsARM.Range("H66").Value = 4749.8
$wsARM.Range("I66").Value = 775
$wsARM.Range("K66").Value = 3875
$wsARM.Range("M66").Value = -443
$wsARM.Range("H74").Value = 2059.9714
$wsARM.Range("J74").Value = 0
$wsARM.Range("L74").Value = 0
$wsARM.Range("N74").ClearContents()
$wsARM.Range("H77").Value = 2059.9714
$wsARM.Range("J77").Value = 0
$wsARM.Range("L77").Value = 0
$wsARM.Range("N77").ClearContents()
$wsARM.Range("H96").Value = 75000
$wsARM.Range("J96").Value = 75000
$wsARM.Range("L96").Value = 75000
$wsARM.Range("N96").Value = -80492
$wsARM.Range("H116").Value = 1544.2693
$wsARM.Range("I116").Value = 571.35297
$wsARM.Range("K116").Value = 571.35297
$wsARM.Range("M116").Value = 1722.64703
$wsARM.Range("H132").Value = 4206.278
$wsARM.Range("J132").Value = 0
$wsARM.Range("L132").Value = 0
$wsARM.Range("N132").ClearContents()
$wsARM.Range("H134").Value = 36499
$wsARM.Range("J134").Value = 36499
$wsARM.Range("L134").Value = 36499
$wsARM.Range("N134").Value = -46639
$wsARM.Range("H136").Value = 3269
$wsARM.Range("I136").Value = 3269
$wsARM.Range("J136").Value = 0
$wsARM.Range("K136").Value = 9807
$wsARM.Range("L136").Value = 0
$wsARM.Range("M136").Value = -7257
$wsARM.Range("N136").ClearContents()

# --- BSM ---
$wsBSM.Range("H3").Value = 1544.2693
$wsBSM.Range("I3").Value = 571.35297
$wsBSM.Range("K3").Value = 571.35297
$wsBSM.Range("M3").Value = -457.35297
$wsBSM.Range("H99").Value = 3402.1428
$wsBSM.Range("I99").Value = 1817.1428
$wsBSM.Range("K99").Value = 1817.1428
$wsBSM.Range("M99").Value = -319.1428000000001

# --- CRP ---
$wsCRP.Range("H31").Value = 2904.054
$wsCRP.Range("I31").Value = 1593.0834
$wsCRP.Range("J31").Value = 5324.3076
$wsCRP.Range("K31").Value = 1593.0834
$wsCRP.Range("L31").Value = 5324.3076
$wsCRP.Range("M31").Value = -1298.0834
$wsCRP.Range("N31").Value = -5914.3076
$wsCRP.Range("H34").Value = 2904.054
$wsCRP.Range("I34").Value = 1593.0834
$wsCRP.Range("J34").Value = 5324.3076
$wsCRP.Range("K34").Value = 1593.0834
$wsCRP.Range("L34").Value = 5324.3076
$wsCRP.Range("M34").Value = -1391.0834
$wsCRP.Range("N34").Value = -5728.3076
$wsCRP.Range("H52").Value = 90709
$wsCRP.Range("J52").Value = 0
$wsCRP.Range("L52").Value = 0
$wsCRP.Range("N52").ClearContents()
$wsCRP.Range("H58").Value = 1683.1666
$wsCRP.Range("I58").Value = 1769.8
$wsCRP.Range("K58").Value = 1769.8
$wsCRP.Range("M58").Value = -1566.8
$wsCRP.Range("H105").Value = 2000
$wsCRP.Range("I105").Value = 2000
$wsCRP.Range("K105").Value = 2000
$wsCRP.Range("M105").Value = -253
$wsCRP.Range("H107").Value = 1470.2858
$wsCRP.Range("I107").Value = 1470.2858
$wsCRP.Range("J107").Value = 0
$wsCRP.Range("K107").Value = 1470.2858
$wsCRP.Range("L107").Value = 0
$wsCRP.Range("M107").Value = 449.7141999999999
$wsCRP.Range("N107").ClearContents()
$wsCRP.Range("H136").Value = 1683.1666
$wsCRP.Range("I136").Value = 1769.8
$wsCRP.Range("K136").Value = 5309.4
$wsCRP.Range("M136").Value = -2759.4

# --- CUL ---
$wsCUL.Range("H107").Value = 3445.16
$wsCUL.Range("I107").Value = 5395.25
$wsCUL.Range("J107").Value = 2527.4707
$wsCUL.Range("K107").Value = 16185.75
$wsCUL.Range("L107").Value = 7582.4121
$wsCUL.Range("M107").Value = -14265.75
$wsCUL.Range("N107").Value = -11422.4121
$wsCUL.Range("H122").Value = 1698.75
$wsCUL.Range("I122").Value = 1900
$wsCUL.Range("K122").Value = 17100
$wsCUL.Range("M122").Value = -14650
$wsCUL.Range("H129").Value = 2421.682
$wsCUL.Range("J129").Value = 4019.4
$wsCUL.Range("L129").Value = 12058.2
$wsCUL.Range("N129").Value = -22058.2
$wsCUL.Range("H131").Value = 2382984
$wsCUL.Range("I131").Value = 778.5454999999999
$wsCUL.Range("J131").Value = 3924411
$wsCUL.Range("K131").Value = 2335.6365
$wsCUL.Range("L131").Value = 11773233
$wsCUL.Range("M131").Value = 2704.3635
$wsCUL.Range("N131").Value = -11783313
$wsCUL.Range("H134").Value = 3093.348
$wsCUL.Range("I134").Value = 1543.1333
$wsCUL.Range("K134").Value = 4629.3999
$wsCUL.Range("M134").Value = 440.6000999999997
$wsCUL.Range("H139").Value = 3895.4736
$wsCUL.Range("I139").Value = 2376.75
$wsCUL.Range("K139").Value = 7130.25
$wsCUL.Range("M139").Value = -1990.25
$wsCUL.Range("H140").Value = 2978.2144
$wsCUL.Range("I140").Value = 2978.2144
$wsCUL.Range("K140").Value = 8934.643199999999
$wsCUL.Range("M140").Value = -3754.643199999999

# --- GSM ---
$wsGSM.Range("H21").Value = 2010000
$wsGSM.Range("I21").Value = 20000
$wsGSM.Range("K21").Value = 20000
$wsGSM.Range("M21").Value = -19827
$wsGSM.Range("H30").Value = 2010000
$wsGSM.Range("I30").Value = 20000
$wsGSM.Range("K30").Value = 20000
$wsGSM.Range("M30").Value = -19895
$wsGSM.Range("H38").Value = 5999.5
$wsGSM.Range("J38").Value = 5999.5
$wsGSM.Range("L38").Value = 5999.5
$wsGSM.Range("N38").Value = -6925.5
$wsGSM.Range("H107").Value = 462.75
$wsGSM.Range("I107").Value = 86.59999999999999
$wsGSM.Range("J107").Value = 1089.6666
$wsGSM.Range("K107").Value = 86.59999999999999
$wsGSM.Range("L107").Value = 1089.6666
$wsGSM.Range("M107").Value = 1833.4
$wsGSM.Range("N107").Value = -4929.6666
$wsGSM.Range("H122").Value = 13691.032
$wsGSM.Range("I122").Value = 15800.92
$wsGSM.Range("K122").Value = 47402.76
$wsGSM.Range("M122").Value = -44952.76

# --- LTW ---
$wsLTW.Range("H7").Value = 3911.1765
$wsLTW.Range("I7").Value = 3974.375
$wsLTW.Range("K7").Value = 3974.375
$wsLTW.Range("M7").Value = -3862.375
$wsLTW.Range("H22").Value = 1429.1136
$wsLTW.Range("I22").Value = 2503.6365
$wsLTW.Range("J22").Value = 1070.9395
$wsLTW.Range("K22").Value = 2503.6365
$wsLTW.Range("L22").Value = 1070.9395
$wsLTW.Range("M22").Value = -2208.6365
$wsLTW.Range("N22").Value = -1660.9395
$wsLTW.Range("H27").Value = 1429.1136
$wsLTW.Range("I27").Value = 2503.6365
$wsLTW.Range("J27").Value = 1070.9395
$wsLTW.Range("K27").Value = 2503.6365
$wsLTW.Range("L27").Value = 1070.9395
$wsLTW.Range("M27").Value = -2396.6365
$wsLTW.Range("N27").Value = -1284.9395
$wsLTW.Range("H46").Value = 3000
$wsLTW.Range("J46").Value = 3000
$wsLTW.Range("L46").Value = 3000
$wsLTW.Range("N46").Value = -3376
$wsLTW.Range("H68").Value = 6215.8066
$wsLTW.Range("I68").Value = 2656.5
$wsLTW.Range("K68").Value = 2656.5
$wsLTW.Range("M68").Value = -1907.5
$wsLTW.Range("H71").Value = 6215.8066
$wsLTW.Range("I71").Value = 2656.5
$wsLTW.Range("K71").Value = 13282.5
$wsLTW.Range("M71").Value = -9538.5
$wsLTW.Range("H99").Value = 0
$wsLTW.Range("J99").Value = 0
$wsLTW.Range("L99").Value = 0
$wsLTW.Range("N99").ClearContents()
$wsLTW.Range("H100").Value = 9582.091
$wsLTW.Range("I100").Value = 8466.666999999999
$wsLTW.Range("K100").Value = 8466.666999999999
$wsLTW.Range("M100").Value = -7925.666999999999
$wsLTW.Range("H126").Value = 3911.1765
$wsLTW.Range("I126").Value = 3974.375
$wsLTW.Range("K126").Value = 11923.125
$wsLTW.Range("M126").Value = -9453.125

# --- WVR ---
$wsWVR.Range("H126").Value = 2726.6667
$wsWVR.Range("I126").Value = 2867.3125
$wsWVR.Range("K126").Value = 8601.9375
$wsWVR.Range("M126").Value = -6131.9375
